# Prix Spot: insert a new date column before column EN (01-oct.) for "11-dec".
$wb = $excel.ActiveWorkbook
$wsPrix = $wb.Worksheets.Item("Prix Spot")

$wsPrix.Range("EN1").EntireColumn.Insert()

$wsPrix.Range("EN1").Value = "11-dec"
for ($r = 2; $r -le 25; $r++) {
    $wsPrix.Cells.Item($r, 144).Value = "-"
}

# Gaz: append the 2025-12-09 row.
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Cells.Item(174, 1).Value = "'2025-12-09"
$wsGaz.Cells.Item(174, 2).Value = 26.1

# CO2: append the 2025-12-09 row.
$wsCO2 = $wb.Worksheets.Item("CO2")
$wsCO2.Cells.Item(174, 1).Value = "'2025-12-09"
$wsCO2.Cells.Item(174, 2).Value = 82.67
